$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh crypto price/volume snapshot data (includes two ranking swaps) as inline text,
# using a leading apostrophe to force text storage for numeric-looking price strings.
$ws.Range('D2').Value = '''26.291.45'
$ws.Range('E2').Value = '''  +0.95%  '
$ws.Range('D3').Value = '''1.679.34'
$ws.Range('D5').Value = '''218.04'
$ws.Range('E5').Value = '''  +0.56%  '
$ws.Range('D6').Value = '''0.5332'
$ws.Range('E6').Value = '''  +4.33%  '
$ws.Range('E7').Value = '''  +0.16%  '
$ws.Range('D8').Value = '''0.2681'
$ws.Range('E8').Value = '''  +1.02%  '
$ws.Range('D9').Value = '''0.06468'
$ws.Range('E9').Value = '''  +1.00%  '
$ws.Range('D10').Value = '''21.93'
$ws.Range('E10').Value = '''  +0.27%  '
$ws.Range('D11').Value = '''0.07528'
$ws.Range('E11').Value = '''  +1.25%  '
$ws.Range('D12').Value = '''1.687.94'
$ws.Range('E12').Value = '''  +0.79%  '
$ws.Range('D13').Value = '''4.526'
$ws.Range('E13').Value = '''  +0.61%  '
$ws.Range('D14').Value = '''0.5777'
$ws.Range('E14').Value = '''  -0.91%  '
$ws.Range('D15').Value = '''0.000008462'
$ws.Range('E15').Value = '''  -1.01%  '
$ws.Range('D16').Value = '''64.75'
$ws.Range('E16').Value = '''  +0.67%  '
$ws.Range('D17').Value = '''26.324.58'
$ws.Range('E17').Value = '''  +0.88%  '
$ws.Range('D18').Value = '''4.897'
$ws.Range('E18').Value = '''  -0.96%  '
$ws.Range('E19').Value = '''  +0.13%  '
$ws.Range('D20').Value = '''10.86'
$ws.Range('E20').Value = '''  +0.95%  '
$ws.Range('D21').Value = '''190.92'
$ws.Range('E21').Value = '''  +0.32%  '
$ws.Range('D22').Value = '''6.207'
$ws.Range('E22').Value = '''  -0.37%  '
$ws.Range('D23').Value = '''1.008'
$ws.Range('E23').Value = '''  +0.14%  '
$ws.Range('D24').Value = '''145.66'
$ws.Range('E24').Value = '''  +0.30%  '
$ws.Range('D25').Value = '''7.823'
$ws.Range('E25').Value = '''  +2.52%  '
$ws.Range('D26').Value = '''0.1273'
$ws.Range('E26').Value = '''  +5.81%  '
$ws.Range('D27').Value = '''15.77'
$ws.Range('E27').Value = '''  +0.84%  '
$ws.Range('D28').Value = '''0.06487'
$ws.Range('E28').Value = '''  -1.29%  '
$ws.Range('D29').Value = '''1.381'
$ws.Range('E29').Value = '''  +4.15%  '
$ws.Range('D30').Value = '''1.322'
$ws.Range('E30').Value = '''  +0.37%  '
$ws.Range('D31').Value = '''3.579'
$ws.Range('E31').Value = '''  +0.98%  '
$ws.Range('D32').Value = '''3.583'
$ws.Range('E32').Value = '''  +1.86%  '
$ws.Range('E33').Value = '''  +1.17%  '
$ws.Range('E34').Value = '''  +1.50%  '
$ws.Range('D35').Value = '''0.6180'
$ws.Range('E35').Value = '''  +1.35%  '
$ws.Range('D36').Value = '''2.401'
$ws.Range('E36').Value = '''  +1.34%  '
$ws.Range('E37').Value = '''  -0.32%  '
$ws.Range('D38').Value = '''6.251'
$ws.Range('E38').Value = '''  +0.22%  '
$ws.Range('D39').Value = '''1.111.48'
$ws.Range('E39').Value = '''  +2.29%  '
$ws.Range('D40').Value = '''0.01620'
$ws.Range('E40').Value = '''  +0.98%  '
$ws.Range('E41').Value = '''  +1.14%  '
$ws.Range('E42').Value = '''  +0.52%  '
$ws.Range('D43').Value = '''100.36'
$ws.Range('E43').Value = '''  -0.26%  '
$ws.Range('D44').Value = '''1.829.73'
$ws.Range('E44').Value = '''  +0.72%  '
$ws.Range('E45').Value = '''  -4.89%  '
$ws.Range('D46').Value = '''57.12'
$ws.Range('E46').Value = '''  +1.40%  '
$ws.Range('B47').Value = '''Frax'
$ws.Range('C47').Value = '''https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D47').Value = '''1.007'
$ws.Range('E47').Value = '''  -0.22%  '
$ws.Range('B48').Value = '''EnergySwap'
$ws.Range('C48').Value = '''https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '''8.157'
$ws.Range('E48').Value = '''  +1.26%  '
$ws.Range('D49').Value = '''0.05264'
$ws.Range('E49').Value = '''  +0.49%  '
$ws.Range('B50').Value = '''Mantle'
$ws.Range('C50').Value = '''https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').Value = '''0.4289'
$ws.Range('E50').Value = '''  +0.01%  '
$ws.Range('B51').Value = '''Aptos'
$ws.Range('C51').Value = '''https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D51').Value = '''6.079'
$ws.Range('E51').Value = '''  +1.35%  '
